$wb = $excel.ActiveWorkbook

# --- Template sheet: insert a new "Group" column ahead of "Collection" (col M) ---
$wsTemplate = $wb.Worksheets.Item("Template")
$wsTemplate.Columns("M:M").Insert()
$wsTemplate.Range("M1").Value = "Group"

# --- Sample Data sheet: same column insert, plus a sample "Bonell" group value ---
$wsSample = $wb.Worksheets.Item("Sample Data")
$wsSample.Columns("M:M").Insert()
$wsSample.Range("M1").Value = "Group"
$wsSample.Range("M2").Value = "Bonell"

# --- View state: Sample Data is no longer the active tab; selection moves to M3 ---
$wsSample.Activate()
[void]$wsSample.Range("M3").Select()

# --- Template becomes the active sheet/tab, with L2 selected ---
$wsTemplate.Activate()
[void]$wsTemplate.Range("L2").Select()
